# Scheduled-runner market-data refresh: rewrite the price/profit columns
# (currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
# LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) for the rows whose
# Universalis market snapshot changed, sheet by sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 287971.06
$ws.Range("I132").Value = 320787.25
$ws.Range("K132").Value = 962361.75
$ws.Range("M132").Value = -959831.75
# Row 137
$ws.Range("H137").Value = 23810764
$ws.Range("I137").Value = 41667636
$ws.Range("J137").Value = 1604.3334
$ws.Range("K137").Value = 125002908
$ws.Range("L137").Value = 4813.0002
$ws.Range("M137").Value = -125000358
$ws.Range("N137").Value = -9913.0002
# Row 138
$ws.Range("H138").Value = 1555.35
$ws.Range("I138").Value = 548.322
$ws.Range("J138").Value = 3004.4878
$ws.Range("K138").Value = 1644.966
$ws.Range("L138").Value = 9013.463400000001
$ws.Range("M138").Value = 3495.034
$ws.Range("N138").Value = -19293.4634
# Row 141
$ws.Range("H141").Value = 1876.2051
$ws.Range("I141").Value = 1091.0154
$ws.Range("J141").Value = 5802.154
$ws.Range("K141").Value = 3273.0462
$ws.Range("L141").Value = 17406.462
$ws.Range("M141").Value = 1906.9538
$ws.Range("N141").Value = -27766.462

$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("N19").ClearContents()
# Row 74
$ws.Range("H74").Value = 4833.1387
$ws.Range("I74").Value = 1370.6129
$ws.Range("J74").Value = 26300.8
$ws.Range("K74").Value = 1370.6129
$ws.Range("L74").Value = 26300.8
$ws.Range("M74").Value = -496.6129000000001
$ws.Range("N74").Value = -28048.8
# Row 77
$ws.Range("H77").Value = 4833.1387
$ws.Range("I77").Value = 1370.6129
$ws.Range("J77").Value = 26300.8
$ws.Range("K77").Value = 6853.0645
$ws.Range("L77").Value = 131504
$ws.Range("M77").Value = -2485.0645
$ws.Range("N77").Value = -140240
# Row 102
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378
# Row 122
$ws.Range("H122").Value = 1847.9697
$ws.Range("I122").Value = 1831.8
$ws.Range("K122").Value = 5495.4
$ws.Range("M122").Value = -3045.4
# Row 132
$ws.Range("H132").Value = 2200.5134
$ws.Range("I132").Value = 1913.5667
$ws.Range("J132").Value = 3430.2856
$ws.Range("K132").Value = 5740.7001
$ws.Range("L132").Value = 10290.8568
$ws.Range("M132").Value = -3210.7001
$ws.Range("N132").Value = -15350.8568

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2605.28
$ws.Range("I20").Value = 2446.6
$ws.Range("K20").Value = 2446.6
$ws.Range("M20").Value = -2199.6
# Row 107
$ws.Range("H107").Value = 730.5769
$ws.Range("I107").Value = 676.2632
$ws.Range("J107").Value = 878
$ws.Range("K107").Value = 676.2632
$ws.Range("L107").Value = 878
$ws.Range("M107").Value = 1243.7368
$ws.Range("N107").Value = -4718
# Row 134
$ws.Range("H134").Value = 2406.2856
$ws.Range("I134").Value = 1426.921
$ws.Range("J134").Value = 4473.8335
$ws.Range("K134").Value = 4280.763
$ws.Range("L134").Value = 13421.5005
$ws.Range("M134").Value = -1745.763
$ws.Range("N134").Value = -18491.5005

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1782.1538
$ws.Range("I58").Value = 658.5143
$ws.Range("J58").Value = 4095.5293
$ws.Range("K58").Value = 658.5143
$ws.Range("L58").Value = 4095.5293
$ws.Range("M58").Value = -455.5143
$ws.Range("N58").Value = -4501.5293
# Row 62
$ws.Range("H62").Value = 29037.375
$ws.Range("I62").Value = 36649.832
$ws.Range("J62").Value = 6200
$ws.Range("K62").Value = 36649.832
$ws.Range("L62").Value = 6200
$ws.Range("M62").Value = -36025.832
$ws.Range("N62").Value = -7448
# Row 65
$ws.Range("H65").Value = 29037.375
$ws.Range("I65").Value = 36649.832
$ws.Range("J65").Value = 6200
$ws.Range("K65").Value = 183249.16
$ws.Range("L65").Value = 31000
$ws.Range("M65").Value = -180129.16
$ws.Range("N65").Value = -37240
# Row 109
$ws.Range("H109").Value = 29733.334
$ws.Range("I109").Value = 29500
$ws.Range("J109").Value = 29850
$ws.Range("K109").Value = 29500
$ws.Range("L109").Value = 29850
$ws.Range("M109").Value = -28460
$ws.Range("N109").Value = -31930
# Row 132
$ws.Range("H132").Value = 1997.9259
$ws.Range("I132").Value = 1538.8043
$ws.Range("K132").Value = 4616.4129
$ws.Range("M132").Value = -2086.4129
# Row 134
$ws.Range("H134").Value = 1839.1
$ws.Range("I134").Value = 1129.6666
$ws.Range("J134").Value = 5859.222
$ws.Range("K134").Value = 3388.9998
$ws.Range("L134").Value = 17577.666
$ws.Range("M134").Value = -853.9998000000001
$ws.Range("N134").Value = -22647.666
# Row 136
$ws.Range("H136").Value = 1782.1538
$ws.Range("I136").Value = 658.5143
$ws.Range("J136").Value = 4095.5293
$ws.Range("K136").Value = 1975.5429
$ws.Range("L136").Value = 12286.5879
$ws.Range("M136").Value = 574.4570999999999
$ws.Range("N136").Value = -17386.5879

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("N36").ClearContents()
# Row 75
$ws.Range("H75").Value = 5676
$ws.Range("I75").Value = 5013
$ws.Range("J75").Value = 6007.5
$ws.Range("K75").Value = 15039
$ws.Range("L75").Value = 18022.5
$ws.Range("M75").Value = -14041
$ws.Range("N75").Value = -20018.5
# Row 78
$ws.Range("H78").Value = 5676
$ws.Range("I78").Value = 5013
$ws.Range("J78").Value = 6007.5
$ws.Range("K78").Value = 45117
$ws.Range("L78").Value = 54067.5
$ws.Range("M78").Value = -40125
$ws.Range("N78").Value = -64051.5
# Row 107
$ws.Range("H107").Value = 1154.1052
$ws.Range("I107").Value = 1476.9166
$ws.Range("J107").Value = 600.7143
$ws.Range("K107").Value = 4430.7498
$ws.Range("L107").Value = 1802.1429
$ws.Range("M107").Value = -2510.7498
$ws.Range("N107").Value = -5642.1429
# Row 113
$ws.Range("H113").Value = 1056.2693
$ws.Range("I113").Value = 683.1667
$ws.Range("J113").Value = 1376.0714
$ws.Range("K113").Value = 2049.5001
$ws.Range("L113").Value = 4128.2142
$ws.Range("M113").Value = 120.4998999999998
$ws.Range("N113").Value = -8468.2142

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 38509
$ws.Range("J6").Value = 38509
$ws.Range("L6").Value = 38509
$ws.Range("N6").Value = -38735
# Row 16
$ws.Range("H16").Value = 38509
$ws.Range("J16").Value = 38509
$ws.Range("L16").Value = 38509
$ws.Range("N16").Value = -39009
# Row 122
$ws.Range("H122").Value = 586153.9
$ws.Range("I122").Value = 1011183.25
$ws.Range("J122").Value = 1738.5
$ws.Range("K122").Value = 3033549.75
$ws.Range("L122").Value = 5215.5
$ws.Range("M122").Value = -3031099.75
$ws.Range("N122").Value = -10115.5
# Row 126
$ws.Range("H126").Value = 3450.2083
$ws.Range("I126").Value = 1687.5
$ws.Range("J126").Value = 3802.75
$ws.Range("K126").Value = 5062.5
$ws.Range("L126").Value = 11408.25
$ws.Range("M126").Value = -2592.5
$ws.Range("N126").Value = -16348.25
# Row 132
$ws.Range("H132").Value = 2978.0527
$ws.Range("I132").Value = 2810.1143
$ws.Range("J132").Value = 4937.3335
$ws.Range("K132").Value = 8430.3429
$ws.Range("L132").Value = 14812.0005
$ws.Range("M132").Value = -5900.3429
$ws.Range("N132").Value = -19872.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").ClearContents()
# Row 103
$ws.Range("H103").Value = 350867.34
$ws.Range("I103").Value = 20000
$ws.Range("J103").Value = 516301
$ws.Range("K103").Value = 20000
$ws.Range("L103").Value = 516301
$ws.Range("M103").Value = -18828
$ws.Range("N103").Value = -518645
# Row 115
$ws.Range("H115").Value = 28666.666
$ws.Range("J115").Value = 28666.666
$ws.Range("L115").Value = 28666.666
$ws.Range("N115").Value = -31800.666
# Row 122
$ws.Range("H122").Value = 68205.92999999999
$ws.Range("I122").Value = 92117.17999999999
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 276351.54
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -273901.54
$ws.Range("N122").Value = -12250
# Row 136
$ws.Range("H136").Value = 9553452
$ws.Range("I136").Value = 10785930
$ws.Range("J136").Value = 1746.25
$ws.Range("K136").Value = 32357790
$ws.Range("L136").Value = 5238.75
$ws.Range("M136").Value = -32355240
$ws.Range("N136").Value = -10338.75
